# Update invalidLoginTest ("Login In Application Bad PW") test data so the
# Object/ObjectType columns pull from the external Objects workbook via
# formulas, matching the pattern already used on validLoginTest.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("validLoginTest")
$ws2 = $wb.Worksheets.Item("invalidLoginTest")

# Row 3: SETTEXT username / name -> Adminrrr
$ws2.Range("C3").Formula = "=[1]Objects!`$B`$2"
$ws2.Range("D3").Formula = "=[1]Objects!`$C`$2"

# Row 4: CLICK password / name
$ws2.Range("C4").Formula = "=[1]Objects!`$B`$3"
$ws2.Range("D4").Formula = "=[1]Objects!`$C`$3"

# Row 5: SETTEXT password / name -> Admin123
$ws2.Range("C5").Formula = "=[1]Objects!`$B`$3"
$ws2.Range("D5").Formula = "=[1]Objects!`$C`$3"

# Row 6: CLICKLISTHOMEPAGE - Object/ObjectType cleared (blank), like sheet1
$ws2.Range("C6").ClearContents()
$ws2.Range("D6").ClearContents()

# Row 7: CLICK loginButton / id
$ws2.Range("C7").Formula = "=[1]Objects!`$B`$4"
$ws2.Range("D7").Formula = "=[1]Objects!`$C`$4"

# Restore the view/selection state: validLoginTest was active with C3
# selected, then the user switched to invalidLoginTest and selected D8,
# leaving it the active tab.
[void]$ws1.Activate()
[void]$ws1.Range("C3").Select()
[void]$ws2.Activate()
[void]$ws2.Range("D8").Select()
